$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume figures to the latest scraped values.
# A leading apostrophe forces Excel to store the value as literal text
# (preserving formats like '59.769.03', trailing zeros such as '24.20',
# and the padded '  -0.34%  ' volume strings) instead of auto-converting
# them to numbers/dates.

$ws.Range("D2").Value = "'59.769.03"
$ws.Range("E2").Value = "'  -0.34%  "
$ws.Range("D3").Value = "'2.360.78"
$ws.Range("E3").Value = "'  -2.22%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'559.31"
$ws.Range("E5").Value = "'  +1.22%  "
$ws.Range("D6").Value = "'133.26"
$ws.Range("E6").Value = "'  -2.69%  "
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E8").Value = "'  -1.66%  "
$ws.Range("E9").Value = "'  -0.39%  "
$ws.Range("D10").Value = "'5.62"
$ws.Range("E10").Value = "'  -1.08%  "
$ws.Range("E11").Value = "'  +0.99%  "
$ws.Range("E12").Value = "'  -3.43%  "
$ws.Range("D13").Value = "'24.20"
$ws.Range("E13").Value = "'  -4.26%  "
$ws.Range("D14").Value = "'2.785.46"
$ws.Range("E14").Value = "'  -2.11%  "
$ws.Range("D15").Value = "'59.754.92"
$ws.Range("E15").Value = "'  -0.27%  "
$ws.Range("E16").Value = "'  -0.31%  "
$ws.Range("D17").Value = "'2.352.90"
$ws.Range("E17").Value = "'  -1.74%  "
$ws.Range("D18").Value = "'11.04"
$ws.Range("E18").Value = "'  -2.48%  "
$ws.Range("E19").Value = "'  +0.93%  "
$ws.Range("D20").Value = "'319.71"
$ws.Range("E20").Value = "'  -2.75%  "
$ws.Range("D21").Value = "'6.66"
$ws.Range("E21").Value = "'  -0.28%  "
$ws.Range("D23").Value = "'64.09"
$ws.Range("E23").Value = "'  -2.63%  "
$ws.Range("E24").Value = "'  -1.12%  "
$ws.Range("E25").Value = "'  -0.09%  "
$ws.Range("D26").Value = "'8.38"
$ws.Range("E26").Value = "'  -2.66%  "
$ws.Range("E27").Value = "'  -1.15%  "
$ws.Range("D28").Value = "'1.81"
$ws.Range("E28").Value = "'  +1.95%  "
$ws.Range("E29").Value = "'  -2.48%  "
$ws.Range("D30").Value = "'170.28"
$ws.Range("E30").Value = "'  +0.72%  "
$ws.Range("D31").Value = "'6.06"
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("E32").Value = "'  +8.59%  "
$ws.Range("D33").Value = "'0.398"
$ws.Range("E33").Value = "'  -1.61%  "
$ws.Range("D34").Value = "'18.08"
$ws.Range("E34").Value = "'  -2.84%  "
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E36").Value = "'  +0.17%  "
$ws.Range("E37").Value = "'  +0.00%  "
$ws.Range("D38").Value = "'4.10"
$ws.Range("E38").Value = "'  -2.10%  "
$ws.Range("D40").Value = "'317.19"
$ws.Range("E40").Value = "'  -1.48%  "
$ws.Range("D41").Value = "'38.60"
$ws.Range("E41").Value = "'  -2.32%  "
$ws.Range("D42").Value = "'144.66"
$ws.Range("E42").Value = "'  +2.91%  "
$ws.Range("E43").Value = "'  -3.88%  "
$ws.Range("E44").Value = "'  -0.63%  "
$ws.Range("D45").Value = "'19.32"
$ws.Range("E45").Value = "'  -1.21%  "
$ws.Range("E46").Value = "'  -1.23%  "
$ws.Range("D47").Value = "'0.565"
$ws.Range("E47").Value = "'  -2.34%  "
$ws.Range("E48").Value = "'  -2.78%  "
$ws.Range("D49").Value = "'11.06"
$ws.Range("E49").Value = "'  +0.18%  "
$ws.Range("D50").Value = "'4.66"
$ws.Range("E50").Value = "'  -0.44%  "
$ws.Range("E51").Value = "'  -1.94%  "
